# #5: fund, bonds, otherbonds, antique done
# Rebuild the "基金受益憑證" (fund) sheet with the extended column layout
# (dealer / property_category / category / date / legislator_name /
#  legislator_id / source_file / index) that the other sheets already use.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# ---- Header row (row 1) -----------------------------------------------
$headers = @{
    2  = "name"
    3  = "owner"
    4  = "dealer"
    5  = "quantity"
    6  = "face_value"
    7  = "currency"
    8  = "total"
    9  = "property_category"
    10 = "category"
    11 = "date"
    12 = "legislator_name"
    13 = "legislator_id"
    14 = "source_file"
    15 = "index"
}
foreach ($col in $headers.Keys) {
    $ws.Cells.Item(1, $col).Value = $headers[$col]
}

# ---- Data rows (rows 2-21) ---------------------------------------------
# columns: A index(no), B name, C owner, D dealer, E quantity, F face_value,
#          G currency, H total, I property_category, J category, K date,
#          L legislator_name, M legislator_id, N source_file, O index
$rows = @(
    @{ A=99;  B="德利資源";      C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=13.251;  F=10; G="歐元";   H=5152;   O=99  },
    @{ A=100; B="瀚亞高科技";    C="李薰楓"; D="保誠證券投資信託公司";     E=14157.8; F=10; G="新臺幣"; H=141578; O=100 },
    @{ A=101; B="元大多福";      C="李薰楓"; D="元大證券投資信託公司";     E=1279.4;  F=10; G="新臺幣"; H=12794;  O=101 },
    @{ A=102; B="元大卓越";      C="李薰楓"; D="元大證券投資信託公司";     E=6187.3;  F=10; G="新臺幣"; H=61873;  O=102 },
    @{ A=103; B="元大亞太成長";  C="李薰楓"; D="元大證券投資信託公司";     E=19931.3; F=10; G="新臺幣"; H=199313; O=103 },
    @{ A=104; B="匯豐太平洋精典";C="李薰楓"; D="國泰世華商業銀行古亭分行"; E=9495.7;  F=10; G="新臺幣"; H=94957;  O=104 },
    @{ A=105; B="富邦精準";      C="李薰楓"; D="台北富邦商業銀行襄陽分行"; E=6000.8;  F=10; G="新臺幣"; H=60008;  O=105 },
    @{ A=106; B="具萊德世界礦業";C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=50.41;   F=10; G="美金";   H=14882;  O=106 },
    @{ A=107; B="安本亞太";      C="李薰楓"; D="國泰世華商業銀行古亭分行"; E=117.007; F=10; G="美金";   H=34542;  O=107 },
    @{ A=108; B="德利資源";      C="李薰楓"; D="國泰世華商業銀行古亭分行"; E=41.629;  F=10; G="歐元";   H=16186;  O=108 },
    @{ A=109; B="天達環球能";    C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=14.003;  F=10; G="美金";   H=4134;   O=109 },
    @{ A=110; B="安泰台灣高股息";C="李薰楓"; D="安泰證券信託投資公司";     E=6259.39; F=10; G="新臺幣"; H=62594;  O=110 },
    @{ A=111; B="富達拉丁美洲";  C="李薰楓"; D="星展（台灣）商業銀行大安分行"; E=108.22; F=10; G="美金"; H=31948;  O=111 },
    @{ A=112; B="施羅德亞洲收";  C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=135.93;  F=10; G="美金";   H=40128;  O=112 },
    @{ A=113; B="德盛新興亞";    C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=41.927;  F=10; G="美金";   H=12377;  O=113 },
    @{ A=115; B="BR環球資配";    C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=42.71;   F=10; G="歐元";   H=16607;  O=115 },
    @{ A=116; B="安本亞太";      C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=34.389;  F=10; G="美金";   H=10152;  O=116 },
    @{ A=117; B="安本亞太";      C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=32.131;  F=10; G="美金";   H=9485;   O=117 },
    @{ A=118; B="BR世界礦業";    C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=42.93;   F=10; G="美金";   H=12673;  O=118 },
    @{ A=119; B="BR世界礦業";    C="李薰楓"; D="台北富邦商業銀行古亭分行"; E=32.68;   F=10; G="美金";   H=9647;   O=119 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = $row.B
    $ws.Cells.Item($r, 3).Value  = $row.C
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    $ws.Cells.Item($r, 11).Value = "2012-04-20"
    $ws.Cells.Item($r, 12).Value = "陳節如"
    $ws.Cells.Item($r, 13).Value = 1709
    $ws.Cells.Item($r, 14).Value = "tmpacad1"
    $ws.Cells.Item($r, 15).Value = $row.O
    $r = $r + 1
}
